$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Copy($ws.Range("C1"))
$ws.Range("C1").Value = "Password"

$ws.Range("A2").Value = "Ayushi123"
$ws.Range("B2").Value = "admin"
$ws.Range("C2").Value = "pokemon@123"
